$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Values")
$ws2 = $wb.Worksheets.Item("Test1")

# "Values" sheet: new chirp-trial inputs for the Posterior block (rows 5-6);
# dependent formulas (B5:C5, B6:C6) recalc automatically.
$ws1.Range("A5").Value = 10
$ws1.Range("A6").Value = 5

# "Test1" sheet: append a new pretension / initial-voltage trial block
$ws2.Range("A8").Value = "DIa 2"
$ws2.Range("A9").Value = "Pretension"
$ws2.Range("B9").Value = 1.83
$ws2.Range("C9").Value = "Initial Voltage"
$ws2.Range("D9").Value = 1.8
$ws2.Range("F9").Value = "Pretension"
$ws2.Range("G9").Value = 3.93
$ws2.Range("H9").Value = "Initial Voltage"
$ws2.Range("I9").Value = 3.91

# Drop the stale row-outline bookkeeping on Test1 (no grouped rows remain)
$ws2.Outline.ShowLevels(0, 0)

# Switch the active tab from Test1 back to Values, restoring each sheet's
# own remembered selection
$ws2.Range("D10").Select()
$ws1.Activate()
$ws1.Range("C7").Select()
